$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览": update a handful of "want to go" counts / ticket prices, then
# insert a new row for the newly-added 南宁·第二届北极光动漫展 event before
# the existing 南宁·万圣漫控嘉年华10 row (which shifts down one row).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 603

$ws1.Range("F4").Value = 486
$ws1.Range("G4").Value = 29.9

$ws1.Range("F5").Value = 499
$ws1.Range("F6").Value = 286
$ws1.Range("F7").Value = 2568
$ws1.Range("F8").Value = 437
$ws1.Range("F9").Value = 6961
$ws1.Range("F11").Value = 439
$ws1.Range("F12").Value = 8

# Insert a brand-new row 13 (pushes the old row 13 down to row 14, carrying
# its values and formatting with it).
$ws1.Rows.Item(13).Insert()

# Restore the bold/bordered "index" style on the new A13 cell by copying it
# from a neighbouring cell that already has it, then set its value.
$ws1.Range("A12").Copy($ws1.Range("A13"))
$ws1.Range("A13").Value = 12

$ws1.Range("B13").Value = "'2024-08-24"
$ws1.Range("B13").Style = "Normal"
$ws1.Range("C13").Value = "南宁·第二届北极光动漫展"
$ws1.Range("D13").Value = "民族大道106号 南宁国际会展中心"
$ws1.Range("E13").Value = "2024.08.24 09:00-08.25 17:00"
$ws1.Range("F13").Value = 47
$ws1.Range("G13").Value = 65
$ws1.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=88276"
$ws1.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202406/mTEwC1GY1717576221099.jpeg"

# The old row 13 (万圣漫控嘉年华10) is now row 14: renumber its index and
# bump its "want to go" count.
$ws1.Range("A14").Value = 13
$ws1.Range("F14").Value = 36

# ---------------------------------------------------------------------------
# Sheet "全部类型": same set of numeric updates (row offsets differ because
# this sheet also contains the 演出 rows interleaved), plus the same new row
# insert before 南宁·万圣漫控嘉年华10 (here at row 17 -> shifts to row 18).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 603

$ws4.Range("F4").Value = 486
$ws4.Range("G4").Value = 29.9

$ws4.Range("F5").Value = 499
$ws4.Range("F6").Value = 286
$ws4.Range("F9").Value = 2568
$ws4.Range("F10").Value = 437
$ws4.Range("F11").Value = 6961
$ws4.Range("F13").Value = 439
$ws4.Range("F14").Value = 8

$ws4.Rows.Item(17).Insert()

$ws4.Range("A16").Copy($ws4.Range("A17"))
$ws4.Range("A17").Value = 16

$ws4.Range("B17").Value = "'2024-08-24"
$ws4.Range("B17").Style = "Normal"
$ws4.Range("C17").Value = "南宁·第二届北极光动漫展"
$ws4.Range("D17").Value = "民族大道106号 南宁国际会展中心"
$ws4.Range("E17").Value = "2024.08.24 09:00-08.25 17:00"
$ws4.Range("F17").Value = 47
$ws4.Range("G17").Value = 65
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=88276"
$ws4.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202406/mTEwC1GY1717576221099.jpeg"

$ws4.Range("A18").Value = 17
$ws4.Range("F18").Value = 36
